# Update NATMI LR-pair edge statistics to reflect the corrected
# "expressing cells" counts (Ligand-expressing cells / Receptor-expressing
# cells go from 1 to 3) and the resulting recomputed expression/specificity
# metrics, per Dr Hou's advice.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 35.995988
$ws.Range("H2").Value = 107.987964
$ws.Range("I2").Value = 0.5613901502831141
$ws.Range("J2").Value = 0.561390150283114
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 117.044563
$ws.Range("N2").Value = 351.133689
$ws.Range("O2").Value = 0.3245365645427815
$ws.Range("P2").Value = 0.3245365645427815
$ws.Range("Q2").Value = 4213.134685213245
$ws.Range("R2").Value = 37918.2121669192
$ws.Range("S2").Value = 0.1821916307410376
$ws.Range("T2").Value = 0.1821916307410376

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 35.995988
$ws.Range("H3").Value = 107.987964
$ws.Range("I3").Value = 0.5613901502831141
$ws.Range("J3").Value = 0.561390150283114
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.281657135515876
$ws.Range("P3").Value = 0.281657135515876
$ws.Range("Q3").Value = 3656.473804890219
$ws.Range("R3").Value = 32908.26424401197
$ws.Range("S3").Value = 0.1581195416355691
$ws.Range("T3").Value = 0.158119541635569

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 35.995988
$ws.Range("H4").Value = 107.987964
$ws.Range("I4").Value = 0.5613901502831141
$ws.Range("J4").Value = 0.561390150283114
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.3938062999413425
$ws.Range("P4").Value = 0.3938062999413425
$ws.Range("Q4").Value = 5112.394604521195
$ws.Range("R4").Value = 46011.55144069075
$ws.Range("S4").Value = 0.2210789779065074
$ws.Range("T4").Value = 0.2210789779065074

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 20.666474
$ws.Range("H5").Value = 61.999422
$ws.Range("I5").Value = 0.3223124461726698
$ws.Range("J5").Value = 0.3223124461726698
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 117.044563
$ws.Range("N5").Value = 351.133689
$ws.Range("O5").Value = 0.3245365645427815
$ws.Range("P5").Value = 0.3245365645427815
$ws.Range("Q5").Value = 2418.898418080862
$ws.Range("R5").Value = 21770.08576272776
$ws.Range("S5").Value = 0.1046021739902584
$ws.Range("T5").Value = 0.1046021739902584

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 20.666474
$ws.Range("H6").Value = 61.999422
$ws.Range("I6").Value = 0.3223124461726698
$ws.Range("J6").Value = 0.3223124461726698
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.281657135515876
$ws.Range("P6").Value = 0.281657135515876
$ws.Range("Q6").Value = 2099.301200468362
$ws.Range("R6").Value = 18893.71080421526
$ws.Range("S6").Value = 0.09078160033010915
$ws.Range("T6").Value = 0.09078160033010914

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 20.666474
$ws.Range("H7").Value = 61.999422
$ws.Range("I7").Value = 0.3223124461726698
$ws.Range("J7").Value = 0.3223124461726698
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.3938062999413425
$ws.Range("P7").Value = 0.3938062999413425
$ws.Range("Q7").Value = 2935.19294906081
$ws.Range("R7").Value = 26416.7365415473
$ws.Range("S7").Value = 0.1269286718523022
$ws.Range("T7").Value = 0.1269286718523022

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 7.456917333333333
$ws.Range("H8").Value = 22.370752
$ws.Range("I8").Value = 0.116297403544216
$ws.Range("J8").Value = 0.116297403544216
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 117.044563
$ws.Range("N8").Value = 351.133689
$ws.Range("O8").Value = 0.3245365645427815
$ws.Range("P8").Value = 0.3245365645427815
$ws.Range("Q8").Value = 872.7916306071253
$ws.Range("R8").Value = 7855.124675464128
$ws.Range("S8").Value = 0.03774275981148537
$ws.Range("T8").Value = 0.03774275981148537

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 7.456917333333333
$ws.Range("H9").Value = 22.370752
$ws.Range("I9").Value = 0.116297403544216
$ws.Range("J9").Value = 0.116297403544216
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 101.5800373333333
$ws.Range("N9").Value = 304.740112
$ws.Range("O9").Value = 0.281657135515876
$ws.Range("P9").Value = 0.281657135515876
$ws.Range("Q9").Value = 757.4739411115804
$ws.Range("R9").Value = 6817.265470004224
$ws.Range("S9").Value = 0.03275599355019777
$ws.Range("T9").Value = 0.03275599355019777

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.456917333333333
$ws.Range("H10").Value = 22.370752
$ws.Range("I10").Value = 0.116297403544216
$ws.Range("J10").Value = 0.116297403544216
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 142.0267893333333
$ws.Range("N10").Value = 426.080368
$ws.Range("O10").Value = 0.3938062999413425
$ws.Range("P10").Value = 0.3938062999413425
$ws.Range("Q10").Value = 1059.082027177415
$ws.Range("R10").Value = 9531.738244596736
$ws.Range("S10").Value = 0.04579865018253289
$ws.Range("T10").Value = 0.04579865018253289
